$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values for rows 2-4 (columns A and B)
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = 2

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 2

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 2

# Row 5 stays the same (A5=0, B5=1) - no change needed

# Delete row 6 entirely, shrinking the used range to A1:B5
$ws.Rows.Item(6).Delete()
